$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.1947194719471947
$ws.Range("C2").Value = 0.5511551155115512
$ws.Range("J2").Value = 0.0132013201320132
$ws.Range("P2").Value = 0.1485148514851485
$ws.Range("S2").Value = 0.0924092409240924
$ws.Range("B3").Value = 0.005524861878453038
$ws.Range("C3").Value = 0.01657458563535912
$ws.Range("J3").Value = 0.01657458563535912
$ws.Range("P3").Value = 0.7569060773480663
$ws.Range("S3").Value = 0.2044198895027624
$ws.Range("J4").Value = 0.02702702702702703
$ws.Range("P4").Value = 0.6486486486486487
$ws.Range("S4").Value = 0.3243243243243243
$ws.Range("B6").Value = 0.04602510460251046
$ws.Range("F6").Value = 0.05857740585774059
$ws.Range("J6").Value = 0.3054393305439331
$ws.Range("O6").Value = 0.01673640167364017
$ws.Range("Q6").Value = 0.1464435146443515
$ws.Range("R6").Value = 0.07531380753138076
$ws.Range("S6").Value = 0.3514644351464435
$ws.Range("B7").Value = 0.09554140127388536
$ws.Range("D7").Value = 0.006369426751592357
$ws.Range("F7").Value = 0.08280254777070063
$ws.Range("J7").Value = 0.08917197452229299
$ws.Range("O7").Value = 0.01273885350318471
$ws.Range("Q7").Value = 0.2038216560509554
$ws.Range("R7").Value = 0.08280254777070063
$ws.Range("S7").Value = 0.4267515923566879
$ws.Range("B8").Value = 0.09153318077803203
$ws.Range("D8").Value = 0.02288329519450801
$ws.Range("E8").Value = 0.002288329519450801
$ws.Range("F8").Value = 0.08237986270022883
$ws.Range("J8").Value = 0.09153318077803203
$ws.Range("O8").Value = 0.01830663615560641
$ws.Range("Q8").Value = 0.1739130434782609
$ws.Range("R8").Value = 0.07780320366132723
$ws.Range("S8").Value = 0.4393592677345537
$ws.Range("B9").Value = 0.1008403361344538
$ws.Range("D9").Value = 0.01680672268907563
$ws.Range("E9").Value = 0.004201680672268907
$ws.Range("F9").Value = 0.05042016806722689
$ws.Range("J9").Value = 0.08403361344537816
$ws.Range("O9").Value = 0.01680672268907563
$ws.Range("Q9").Value = 0.2016806722689076
$ws.Range("R9").Value = 0.1050420168067227
$ws.Range("S9").Value = 0.4201680672268908
$ws.Range("B10").Value = 0.1174628034455756
$ws.Range("D10").Value = 0.01722787783868442
$ws.Range("F10").Value = 0.06108065779169929
$ws.Range("J10").Value = 0.1096319498825372
$ws.Range("O10").Value = 0.01409553641346907
$ws.Range("Q10").Value = 0.2255285826155051
$ws.Range("R10").Value = 0.07987470634299139
$ws.Range("S10").Value = 0.375097885669538
$ws.Range("G11").Value = 0.1234042553191489
$ws.Range("J11").Value = 0.1021276595744681
$ws.Range("K11").Value = 0.1957446808510638
$ws.Range("L11").Value = 0.5574468085106383
$ws.Range("S11").Value = 0.02127659574468085
$ws.Range("G12").Value = 0.7194244604316546
$ws.Range("J12").Value = 0.2014388489208633
$ws.Range("K12").Value = 0.02158273381294964
$ws.Range("L12").Value = 0.03597122302158273
$ws.Range("S12").Value = 0.02158273381294964
$ws.Range("G13").Value = 0.625
$ws.Range("J13").Value = 0.3392857142857143
$ws.Range("S13").Value = 0.03571428571428571
$ws.Range("F15").Value = 0.01739130434782609
$ws.Range("H15").Value = 0.1956521739130435
$ws.Range("I15").Value = 0.0782608695652174
$ws.Range("J15").Value = 0.3173913043478261
$ws.Range("K15").Value = 0.04347826086956522
$ws.Range("M15").Value = 0.01739130434782609
$ws.Range("O15").Value = 0.0782608695652174
$ws.Range("S15").Value = 0.2521739130434782
$ws.Range("F16").Value = 0.02551020408163265
$ws.Range("H16").Value = 0.1275510204081633
$ws.Range("I16").Value = 0.1020408163265306
$ws.Range("J16").Value = 0.4591836734693878
$ws.Range("K16").Value = 0.07653061224489796
$ws.Range("M16").Value = 0.02040816326530612
$ws.Range("O16").Value = 0.06122448979591837
$ws.Range("S16").Value = 0.1275510204081633
$ws.Range("F17").Value = 0.02489626556016597
$ws.Range("H17").Value = 0.1556016597510373
$ws.Range("I17").Value = 0.1037344398340249
$ws.Range("J17").Value = 0.4149377593360996
$ws.Range("K17").Value = 0.08713692946058091
$ws.Range("M17").Value = 0.02282157676348548
$ws.Range("N17").Value = 0.002074688796680498
$ws.Range("O17").Value = 0.07053941908713693
$ws.Range("S17").Value = 0.1182572614107884
$ws.Range("F18").Value = 0.03141361256544502
$ws.Range("H18").Value = 0.2094240837696335
$ws.Range("I18").Value = 0.1047120418848168
$ws.Range("J18").Value = 0.4083769633507853
$ws.Range("K18").Value = 0.05759162303664921
$ws.Range("M18").Value = 0.01047120418848168
$ws.Range("O18").Value = 0.06282722513089005
$ws.Range("S18").Value = 0.1151832460732984
$ws.Range("F19").Value = 0.02484939759036145
$ws.Range("H19").Value = 0.1890060240963855
$ws.Range("I19").Value = 0.09789156626506024
$ws.Range("J19").Value = 0.3719879518072289
$ws.Range("K19").Value = 0.08283132530120482
$ws.Range("M19").Value = 0.0286144578313253
$ws.Range("N19").Value = 0.0007530120481927711
$ws.Range("O19").Value = 0.07228915662650602
$ws.Range("S19").Value = 0.1317771084337349
